$d = $word.ActiveDocument

# The flowchart figure is the last inline picture in the document. The
# paragraph immediately following it is currently empty and needs the
# "Figure: ..." caption added (centered, size 14/28-half-points, underlined).
$shapeCount = $d.InlineShapes.Count
$shp = $d.InlineShapes.Item($shapeCount)
$imgParaEnd = $shp.Range.Paragraphs.Item(1).Range.End

$count = $d.Paragraphs.Count
$targetIndex = -1
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Start -eq $imgParaEnd) {
        $targetIndex = $i
        break
    }
}

$targetPara = $d.Paragraphs.Item($targetIndex)

# Insert the caption text into the (currently empty) paragraph.
$targetPara.Range.InsertBefore("Figure: Flowchart of the General Problem Solution Approach")

# Center the paragraph and underline the text (including the paragraph mark).
$targetPara.Alignment = 1
$targetPara.Range.Font.Underline = 1
$targetPara.Range.Font.Size = 14
